$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the confusion-matrix counts on row 2 (VP, VN, FP, FN)
$ws.Range("B2").Value = 1177
$ws.Range("C2").Value = 1180
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 5

# Clear the center-aligned/bordered style previously on C2:E2 so they
# match the plain (unstyled) look of B2
$ws.Range("C2:E2").Style = "Normal"

# New value entered near the "Total de QRS" block
$ws.Range("G16").Value = 1182
$ws.Range("G16").Style = "Normal"

# Move the active selection to G16 (was J18)
$ws.Range("G16").Select()
